# Apply "Penalty Reward System" forecast/summary update
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet 1: Forecast Comparison ---
# Week_Start_Date (col B) shifts forward by one week; MyForecast (col D) updated.
$weekDates = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02", "2025-02-09",
    "2025-02-16", "2025-02-23", "2025-03-02", "2025-03-09", "2025-03-16",
    "2025-03-23", "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20",
    "2025-04-27"
)
$myForecast = @(43, 43, 43, 43, 44, 45, 45, 45, 46, 46, 41, 42, 41, 41, 41, 41)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    $dateCell = $ws1.Cells.Item($row, 2)
    # Force the date-looking string to stay as literal text (avoid Excel's
    # automatic conversion of "yyyy-mm-dd" strings into date serial values).
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekDates[$i]
    $dateCell.Style = "Normal"

    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# --- Sheet 2: Summary ---
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws2.Cells.Item(2, 2) "2024-07-14 to 2025-01-05"
Set-TextValue $ws2.Cells.Item(5, 2) "21"
Set-TextValue $ws2.Cells.Item(6, 2) "14"
Set-TextValue $ws2.Cells.Item(8, 2) "558 units"
Set-TextValue $ws2.Cells.Item(9, 2) "692"
Set-TextValue $ws2.Cells.Item(10, 2) "351"
Set-TextValue $ws2.Cells.Item(11, 2) "172"
Set-TextValue $ws2.Cells.Item(12, 2) "46"
Set-TextValue $ws2.Cells.Item(14, 2) "41"
Set-TextValue $ws2.Cells.Item(15, 2) "2025-04-27"
